$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Insert a new product row at row 51 ("راجون اخضر 250مل"), pushing the
# existing rows 51-59 (items 45-51, the totals row and the footer row) down
# by one. Excel shifts merged ranges and row content automatically; only the
# brand-new row needs its formatting/values/merges applied explicitly.
# ---------------------------------------------------------------------------
$ws.Rows("51:51").Insert(-4121, 0)

# --- restore the thin light-grey bottom border used by every data row -----
# (mirrors borderId used by styles 7-12 in the original sheet; re-applying
# the identical line style/weight/color lets the engine de-duplicate onto
# the existing style records instead of synthesising new ones)
$newRowBorder = $ws.Range("A51:Q51").Borders.Item(9)
$newRowBorder.LineStyle = 1
$newRowBorder.Weight = 2
$newRowBorder.Color = 13882323

# --- fill in the new row's values ------------------------------------------
$ws.Range("A51").Value = 45
$ws.Range("C51").Value = "راجون اخضر 250مل"
$ws.Range("H51").Value = "1:0"

# L51 and P51 carry numeric-looking text in cells whose number format is
# numeric/date, so flip the format to Text while assigning, then restore the
# original number format, keeping the display intact without reconverting
# the cell to an actual number.
$ws.Range("L51").NumberFormat = "@"
$ws.Range("L51").Value = "0"
$ws.Range("L51").NumberFormat = "#,##0.##;""[""#,##0.##""]"";0"

$ws.Range("N51").Value = "35.00"

$ws.Range("P51").NumberFormat = "@"
$ws.Range("P51").Value = "35.0000"
$ws.Range("P51").NumberFormat = "0.00"

$ws.Range("Q51").Value = "1:0"

# --- merge the new row's cells the same way every other data row is merged
$ws.Range("A51:B51").Merge()
$ws.Range("C51:G51").Merge()
$ws.Range("H51:K51").Merge()
$ws.Range("L51:M51").Merge()
$ws.Range("N51:O51").Merge()

# --- renumber the items that got pushed down (45-51 -> 46-52) -------------
$ws.Range("A52").Value = 46
$ws.Range("A53").Value = 47
$ws.Range("A54").Value = 48
$ws.Range("A55").Value = 49
$ws.Range("A56").Value = 50
$ws.Range("A57").Value = 51
$ws.Range("A58").Value = 52

# --- row heights follow the report's own autofit cycle, keyed off the
# absolute row number; re-assert them for every row from the new one down
# through the totals row (the footer row keeps its fixed 16.5pt height).
$ws.Rows("51:51").RowHeight = 25.5
$ws.Rows("52:52").RowHeight = 25.5
$ws.Rows("53:53").RowHeight = 24.75
$ws.Rows("54:54").RowHeight = 25.5
$ws.Rows("55:55").RowHeight = 24.75
$ws.Rows("56:56").RowHeight = 25.5
$ws.Rows("57:57").RowHeight = 25.5
$ws.Rows("58:58").RowHeight = 24.75
$ws.Rows("59:59").RowHeight = 25.5

# --- update the running total (was P58, now P59) ---------------------------
$ws.Range("P59").Value = 2552.7 + 35

# --- refresh the export timestamp in the footer (was row 59, now row 60) --
$ws.Range("A60").Value = "Friday, 22 August, 2025 7:51 PM"
